# Update "想去人数" (column F) figures across all sheets to reflect the
# latest scrape output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$updates = @{
    3  = 103
    5  = 8135
    7  = 83
    8  = 86
    9  = 7063
    10 = 1133
    11 = 539
    14 = 703
    18 = 230
    21 = 62
    22 = 11574
    23 = 4
    24 = 126
    25 = 2250
    27 = 3135
    29 = 2675
    31 = 23
    32 = 284
    33 = 44
    35 = 1608
    36 = 73
    37 = 99
    38 = 5792
    40 = 1784
    42 = 836
    46 = 1102
    47 = 1514
    48 = 98
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet "演出" -----------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$updates = @{
    13 = 10
    22 = 3
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet "本地生活" --------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$updates = @{
    2 = 231
    3 = 368
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet "全部类型" --------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    3  = 231
    4  = 368
    7  = 8136
    8  = 83
    10 = 86
    11 = 7063
    12 = 7063
    13 = 1133
    14 = 539
    16 = 703
    20 = 230
    22 = 62
    25 = 11574
    27 = 4
    28 = 126
    29 = 2251
    30 = 2251
    31 = 3135
    32 = 2675
    33 = 23
    34 = 284
    35 = 44
    38 = 1608
    39 = 73
    40 = 99
    41 = 5792
    43 = 1784
    46 = 836
    49 = 1102
    50 = 1514
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
